$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 835
$ws.Range("E3").Value = 1573
$ws.Range("E4").Value = 493
$ws.Range("E5").Value = 1907
$ws.Range("E6").Value = 979
$ws.Range("E7").Value = 1186
$ws.Range("E8").Value = 2107
$ws.Range("E9").Value = 2298
$ws.Range("E10").Value = 2031
$ws.Range("E11").Value = 2288
$ws.Range("E12").Value = 2393
$ws.Range("E13").Value = 1937
$ws.Range("E14").Value = 1476
$ws.Range("E15").Value = 1483
$ws.Range("E16").Value = 2413

$wb.Save()
